# Updates the cryptos price/volume table with freshly scraped values.
# Columns: A=rank(idx), B=Coin, C=Link, D=Price, E=Volume(1h)
# All of these cells store plain text (inline strings) in the source file,
# so numeric-looking values are written with a leading apostrophe to force
# Excel to keep them as text instead of silently re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "28.708.53"
$ws.Cells.Item(2, 5).Value = "  -2.54%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "1.851.78"
$ws.Cells.Item(3, 5).Value = "  -3.41%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).Value = "'1.004"
$ws.Cells.Item(4, 5).Value = "  -0.81%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "'335.49"
$ws.Cells.Item(5, 5).Value = "  +3.08%  "

# Row 6 - USDC
$ws.Cells.Item(6, 4).Value = "'1.004"
$ws.Cells.Item(6, 5).Value = "  -0.73%  "

# Row 7 - XRP
$ws.Cells.Item(7, 4).Value = "'0.4643"
$ws.Cells.Item(7, 5).Value = "  -3.45%  "

# Row 8 - Cardano
$ws.Cells.Item(8, 4).Value = "'0.3929"
$ws.Cells.Item(8, 5).Value = "  -2.93%  "

# Row 9 - OKB
$ws.Cells.Item(9, 4).Value = "'46.65"
$ws.Cells.Item(9, 5).Value = "  -2.49%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 4).Value = "'0.07916"
$ws.Cells.Item(10, 5).Value = "  -3.56%  "

# Row 11 - Polygon
$ws.Cells.Item(11, 4).Value = "'0.9832"
$ws.Cells.Item(11, 5).Value = "  -2.28%  "

# Row 12 - Solana (price unchanged)
$ws.Cells.Item(12, 5).Value = "  -4.51%  "

# Row 13 - WrappedEther
$ws.Cells.Item(13, 4).Value = "1.855.99"
$ws.Cells.Item(13, 5).Value = "  -2.43%  "

# Row 14 - Polkadot: unchanged

# Row 15 - Chainlink
$ws.Cells.Item(15, 4).Value = "'7.012"
$ws.Cells.Item(15, 5).Value = "  -2.90%  "

# Row 16 - TRON
$ws.Cells.Item(16, 4).Value = "'0.06767"
$ws.Cells.Item(16, 5).Value = "  -1.35%  "

# Row 17 - BinanceUSD (price unchanged)
$ws.Cells.Item(17, 5).Value = "  -0.70%  "

# Row 18 - Litecoin
$ws.Cells.Item(18, 4).Value = "'87.60"
$ws.Cells.Item(18, 5).Value = "  -3.94%  "

# Row 19 - ShibaInu
$ws.Cells.Item(19, 4).Value = "'0.00001014"
$ws.Cells.Item(19, 5).Value = "  -2.34%  "

# Row 20 - Avalanche
$ws.Cells.Item(20, 4).Value = "'17.05"
$ws.Cells.Item(20, 5).Value = "  -2.62%  "

# Row 21 - Dai
$ws.Cells.Item(21, 4).Value = "'1.004"
$ws.Cells.Item(21, 5).Value = "  -0.68%  "

# Row 22 - WrappedBTC
$ws.Cells.Item(22, 4).Value = "28.720.78"
$ws.Cells.Item(22, 5).Value = "  -2.55%  "

# Row 23 - Uniswap
$ws.Cells.Item(23, 4).Value = "'5.416"
$ws.Cells.Item(23, 5).Value = "  -4.33%  "

# Row 24 - Cosmos
$ws.Cells.Item(24, 4).Value = "'11.31"
$ws.Cells.Item(24, 5).Value = "  -4.50%  "

# Row 25 - Toncoin
$ws.Cells.Item(25, 4).Value = "'2.133"
$ws.Cells.Item(25, 5).Value = "  -2.70%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Cells.Item(26, 4).Value = "2.068.49"
$ws.Cells.Item(26, 5).Value = "  -3.03%  "

# Row 27 - Monero
$ws.Cells.Item(27, 4).Value = "'153.46"
$ws.Cells.Item(27, 5).Value = "  -1.51%  "

# Row 28 - InternetComputer(DFINITY) (price unchanged)
$ws.Cells.Item(28, 5).Value = "  -4.49%  "

# Row 29 - EthereumClassic
$ws.Cells.Item(29, 4).Value = "'19.41"
$ws.Cells.Item(29, 5).Value = "  -2.94%  "

# Row 30 - LidoDAOToken
$ws.Cells.Item(30, 4).Value = "'2.027"
$ws.Cells.Item(30, 5).Value = "  -3.24%  "

# Row 31 - BitcoinCash
$ws.Cells.Item(31, 4).Value = "'117.15"
$ws.Cells.Item(31, 5).Value = "  -2.73%  "

# Row 32 - ImmutableX
$ws.Cells.Item(32, 4).Value = "'0.9805"
$ws.Cells.Item(32, 5).Value = "  -3.07%  "

# Row 33 - Stellar
$ws.Cells.Item(33, 4).Value = "'0.09425"
$ws.Cells.Item(33, 5).Value = "  -1.68%  "

# Row 34 - Filecoin
$ws.Cells.Item(34, 4).Value = "'5.397"
$ws.Cells.Item(34, 5).Value = "  -3.68%  "

# Row 35 - HuobiToken (volume unchanged)
$ws.Cells.Item(35, 4).Value = "'3.490"

# Row 36 - ARBITRUM
$ws.Cells.Item(36, 4).Value = "'1.352"
$ws.Cells.Item(36, 5).Value = "  -1.30%  "

# Row 37 - Hedera
$ws.Cells.Item(37, 4).Value = "'0.06132"
$ws.Cells.Item(37, 5).Value = "  -2.32%  "

# Row 38 - VeChain
$ws.Cells.Item(38, 4).Value = "'0.02201"
$ws.Cells.Item(38, 5).Value = "  -3.54%  "

# Row 39 - TrustWalletToken
$ws.Cells.Item(39, 4).Value = "'1.166"
$ws.Cells.Item(39, 5).Value = "  -1.19%  "

# Row 40 - TheSandbox
$ws.Cells.Item(40, 4).Value = "'0.5730"
$ws.Cells.Item(40, 5).Value = "  -3.26%  "

# Row 41 - FraxShare
$ws.Cells.Item(41, 4).Value = "'7.635"
$ws.Cells.Item(41, 5).Value = "  -2.86%  "

# Row 42 - Aptos
$ws.Cells.Item(42, 4).Value = "'10.13"
$ws.Cells.Item(42, 5).Value = "  -5.35%  "

# Row 43 - Algorand (price unchanged)
$ws.Cells.Item(43, 5).Value = "  -2.95%  "

# Row 44 - RenderToken
$ws.Cells.Item(44, 4).Value = "'2.390"
$ws.Cells.Item(44, 5).Value = "  +0.00%  "

# Row 45 - WEMIXToken
$ws.Cells.Item(45, 4).Value = "'1.226"
$ws.Cells.Item(45, 5).Value = "  -4.28%  "

# Rows 46 & 47 swap places: EnergySwap moves to row 46, Decentraland to row 47,
# each also getting freshly scraped Price/Volume values.
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "'11.89"
$ws.Cells.Item(46, 5).Value = "  -4.34%  "

$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47, 4).Value = "'0.5409"
$ws.Cells.Item(47, 5).Value = "  -2.65%  "

# Row 48 - Cronos (price unchanged)
$ws.Cells.Item(48, 5).Value = "  -4.42%  "

# Row 49 - NEARProtocol
$ws.Cells.Item(49, 4).Value = "'1.924"
$ws.Cells.Item(49, 5).Value = "  -0.27%  "

# Row 50 - Quant
$ws.Cells.Item(50, 4).Value = "'115.91"
$ws.Cells.Item(50, 5).Value = "  -1.75%  "

# Row 51 - Elrond
$ws.Cells.Item(51, 4).Value = "'43.59"
$ws.Cells.Item(51, 5).Value = "  +4.11%  "
